# Daily Status Tracker Updated
#
# 1. Row 4 (SlNo 3) status flips from "WIP" to "Done".
# 2. Row 18 (SlNo 17) comment moves from "To be finished by 15-Feb" to
#    "To be finished by 19-Feb".
# 3. The sheet is filtered on "Owned by" (column D) to show only "Rahul",
#    which hides every other row.
# 4. The active selection ends up on G18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------------
$ws.Range("F4").Value = "Done"
$ws.Range("G18").Value = "To be finished by 19-Feb"

# --- Filter "Owned by" to Rahul only (xlFilterValues-style checkbox list) -
$ws.Range("A1:H23").AutoFilter(4, @("Rahul"), 7)

# --- Selection matches the post-edit workbook ------------------------------
$ws.Range("G18").Select()
